$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7760525136293381
$ws.Range("C2").Value = 0.2180509122735828
$ws.Range("E2").Value = 0.6142252801313077
$ws.Range("F2").Value = 1.757488351347746
$ws.Range("G2").Value = 0.2263892630966495
$ws.Range("H2").Value = 0.4108106610496733
$ws.Range("J2").Value = 0.02531087551720645
$ws.Range("M2").Value = 0.5844291644989994
$ws.Range("O2").Value = 1.188163707121475
$ws.Range("B3").Value = 0.6778198579905279
$ws.Range("C3").Value = 0.1975436041533385
$ws.Range("E3").Value = 0.6094916102050902
$ws.Range("F3").Value = 1.754147175247581
$ws.Range("G3").Value = 0.229686122775064
$ws.Range("H3").Value = 0.4165103899592637
$ws.Range("J3").Value = 0.02523304911029456
$ws.Range("M3").Value = 0.5423405840568734
$ws.Range("O3").Value = 1.20682936809321
$ws.Range("B4").Value = 0.6172785654531481
$ws.Range("C4").Value = 0.184912588078646
$ws.Range("E4").Value = 0.6068645613842847
$ws.Range("F4").Value = 1.753280031083364
$ws.Range("G4").Value = 0.2320194039267562
$ws.Range("H4").Value = 0.4202886834045643
$ws.Range("J4").Value = 0.02519748287690327
$ws.Range("M4").Value = 0.5165698385727993
$ws.Range("O4").Value = 1.219519924096545
$ws.Range("B5").Value = 0.5925523104570232
$ws.Range("C5").Value = 0.1797558400093351
$ws.Range("E5").Value = 0.6058643311945247
$ws.Range("F5").Value = 1.753224164932519
$ws.Range("G5").Value = 0.2330476092990423
$ws.Range("H5").Value = 0.4218983654264008
$ws.Range("J5").Value = 0.02518607236251569
$ws.Range("M5").Value = 0.5060867275251013
$ws.Range("O5").Value = 1.224999824122577
$ws.Range("B6").Value = 0.5884432504918777
$ws.Range("C6").Value = 0.1788990039769658
$ws.Range("E6").Value = 0.6057024918553751
$ws.Range("F6").Value = 1.753232847668144
$ws.Range("G6").Value = 0.2332230056826248
$ws.Range("H6").Value = 0.4221698775141931
$ws.Range("J6").Value = 0.025184364168096
$ws.Range("M6").Value = 0.5043471612833557
$ws.Range("O6").Value = 1.22592835902374
$ws.Range("B7").Value = 0.6169453197763346
$ws.Range("C7").Value = 0.1848430803117651
$ws.Range("E7").Value = 0.6068507871736912
$ws.Range("F7").Value = 1.753278073495466
$ws.Range("G7").Value = 0.2320329578101479
$ws.Range("H7").Value = 0.420310108816949
$ws.Range("J7").Value = 0.02519731649408286
$ws.Range("M7").Value = 0.5164283832248344
$ws.Range("O7").Value = 1.21959258037954
$ws.Range("B8").Value = 0.7422297678719474
$ws.Range("C8").Value = 0.2109884193797598
$ws.Range("E8").Value = 0.6125351361599343
$ws.Range("F8").Value = 1.756090307155645
$ws.Range("G8").Value = 0.2274616929812368
$ws.Range("H8").Value = 0.4127180555936434
$ws.Range("J8").Value = 0.02528151087322783
$ws.Range("M8").Value = 0.5699024738093925
$ws.Range("O8").Value = 1.194343886478592
$ws.Range("B9").Value = 0.9860580427597938
$ws.Range("C9").Value = 0.2619308300883745
$ws.Range("E9").Value = 0.6258987158461764
$ws.Range("F9").Value = 1.771019546550406
$ws.Range("G9").Value = 0.2209638730881309
$ws.Range("H9").Value = 0.4000439312112007
$ws.Range("J9").Value = 0.02554316318042638
$ws.Range("M9").Value = 0.6753131117226587
$ws.Range("O9").Value = 1.15462657014281
$ws.Range("B10").Value = 1.164004701926046
$ws.Range("C10").Value = 0.2991408766451968
$ws.Range("E10").Value = 0.6370692740184225
$ws.Range("F10").Value = 1.787755533121739
$ws.Range("G10").Value = 0.2177138278962332
$ws.Range("H10").Value = 0.3920859253160174
$ws.Range("J10").Value = 0.02579377319429099
$ws.Range("M10").Value = 0.7530713973138603
$ws.Range("O10").Value = 1.131469607009677
$ws.Range("B11").Value = 1.244685801467824
$ws.Range("C11").Value = 0.3160181771187354
$ws.Range("E11").Value = 0.6424450627900384
$ws.Range("F11").Value = 1.796627856655121
$ws.Range("G11").Value = 0.2165705416260906
$ws.Range("H11").Value = 0.3887604261205979
$ws.Range("J11").Value = 0.02592036164913836
$ws.Range("M11").Value = 0.7885097375974226
$ws.Range("O11").Value = 1.122253914653299
$ws.Range("B12").Value = 1.275197748310859
$ws.Range("C12").Value = 0.322401667456063
$ws.Range("E12").Value = 0.6445230385680958
$ws.Range("F12").Value = 1.800169050410744
$ws.Range("G12").Value = 0.2161861556577378
$ws.Range("H12").Value = 0.3875435915012062
$ws.Range("J12").Value = 0.02597009822203233
$ws.Range("M12").Value = 0.8019382631061376
$ws.Range("O12").Value = 1.118954680706196
$ws.Range("B13").Value = 1.268628271530815
$ws.Range("C13").Value = 0.3210272123044433
$ws.Range("E13").Value = 0.6440736293189531
$ws.Range("F13").Value = 1.799398315314463
$ws.Range("G13").Value = 0.2162667750867513
$ws.Range("H13").Value = 0.387803768749535
$ws.Range("J13").Value = 0.02595930664136503
$ws.Range("M13").Value = 0.7990458090562527
$ws.Range("O13").Value = 1.119656740040455
$ws.Range("B14").Value = 1.247196853940864
$ws.Range("C14").Value = 0.3165435049824907
$ws.Range("E14").Value = 0.6426151719032731
$ws.Range("F14").Value = 1.796915554296504
$ws.Range("G14").Value = 0.2165379428916623
$ws.Range("H14").Value = 0.3886594649262847
$ws.Range("J14").Value = 0.02592441748030794
$ws.Range("M14").Value = 0.7896143379303595
$ws.Range("O14").Value = 1.121978660405929
$ws.Range("B15").Value = 1.234064200894693
$ws.Range("C15").Value = 0.3137961054940206
$ws.Range("E15").Value = 0.6417273299849171
$ws.Range("F15").Value = 1.795418431516907
$ws.Range("G15").Value = 0.2167103743421421
$ws.Range("H15").Value = 0.3891891357580874
$ws.Range("J15").Value = 0.02590328103375938
$ws.Range("M15").Value = 0.7838384156170122
$ws.Range("O15").Value = 1.123425746752801
$ws.Range("B16").Value = 1.158726466212102
$ws.Range("C16").Value = 0.2980368681004961
$ws.Range("E16").Value = 0.6367238718935013
$ws.Range("F16").Value = 1.78720107428633
$ws.Range("G16").Value = 0.2177953192255728
$ws.Range("H16").Value = 0.3923091911751513
$ws.Range("J16").Value = 0.02578575277116002
$ws.Range("M16").Value = 0.7507566846548457
$ws.Range("O16").Value = 1.132098487289923
$ws.Range("B17").Value = 1.112439394421074
$ws.Range("C17").Value = 0.2883560437700794
$ws.Range("E17").Value = 0.6337297535536663
$ws.Range("F17").Value = 1.782482739141557
$ws.Range("G17").Value = 0.2185469966599598
$ws.Range("H17").Value = 0.3942987728629177
$ws.Range("J17").Value = 0.02571686950113872
$ws.Range("M17").Value = 0.7304784768163159
$ws.Range("O17").Value = 1.137757313623212
$ws.Range("B18").Value = 1.085791197813933
$ws.Range("C18").Value = 0.2827832333448441
$ws.Range("E18").Value = 0.6320353134834846
$ws.Range("F18").Value = 1.779887361104059
$ws.Range("G18").Value = 0.2190108689114467
$ws.Range("H18").Value = 0.3954708548728689
$ws.Range("J18").Value = 0.02567843474960085
$ws.Range("M18").Value = 0.7188212113353813
$ws.Range("O18").Value = 1.141136170504723
$ws.Range("B19").Value = 1.076764319013705
$ws.Range("C19").Value = 0.280895588791708
$ws.Range("E19").Value = 0.631466363238701
$ws.Range("F19").Value = 1.779028948489724
$ws.Range("G19").Value = 0.2191733321127671
$ws.Range("H19").Value = 0.3958724611009998
$ws.Range("J19").Value = 0.02566562523112381
$ws.Range("M19").Value = 0.7148753503892635
$ws.Range("O19").Value = 1.142301470133162
$ws.Range("B20").Value = 1.117369342317772
$ws.Range("C20").Value = 0.2893870690539586
$ws.Range("E20").Value = 0.6340456161899581
$ws.Range("F20").Value = 1.782972748546214
$ws.Range("G20").Value = 0.218463713789923
$ws.Range("H20").Value = 0.3940841081674051
$ws.Range("J20").Value = 0.02572407965961077
$ws.Range("M20").Value = 0.7326364873150482
$ws.Range("O20").Value = 1.137142076355801
$ws.Range("B21").Value = 1.253492886496304
$ws.Range("C21").Value = 0.317860687512848
$ws.Range("E21").Value = 0.643042408941696
$ws.Range("F21").Value = 1.797639874136536
$ws.Range("G21").Value = 0.2164569738426607
$ws.Range("H21").Value = 0.3884069731451518
$ws.Range("J21").Value = 0.02593461648783091
$ws.Range("M21").Value = 0.7923843561076609
$ws.Range("O21").Value = 1.121291477153108
$ws.Range("B22").Value = 1.342222180128317
$ws.Range("C22").Value = 0.3364254843910999
$ws.Range("E22").Value = 0.6491687667410346
$ws.Range("F22").Value = 1.808283424352936
$ws.Range("G22").Value = 0.2154286128558951
$ws.Range("H22").Value = 0.3849441673162275
$ws.Range("J22").Value = 0.0260827035850788
$ws.Range("M22").Value = 0.831484007297405
$ws.Range("O22").Value = 1.112043294166568
$ws.Range("B23").Value = 1.294887786441905
$ws.Range("C23").Value = 0.3265213009313754
$ws.Range("E23").Value = 0.6458764756038278
$ws.Range("F23").Value = 1.802505850122415
$ws.Range("G23").Value = 0.2159514428932141
$ws.Range("H23").Value = 0.3867696528621423
$ws.Range("J23").Value = 0.02600271002469867
$ws.Range("M23").Value = 0.8106113389470124
$ws.Range("O23").Value = 1.116877246261026
$ws.Range("B24").Value = 1.115140628438667
$ws.Range("C24").Value = 0.2889209648427595
$ws.Range("E24").Value = 0.6339027308399068
$ws.Range("F24").Value = 1.782750850077875
$ws.Range("G24").Value = 0.2185012671876194
$ws.Range("H24").Value = 0.3941810701006077
$ws.Range("J24").Value = 0.02572081631075918
$ws.Range("M24").Value = 0.7316608477729147
$ws.Range("O24").Value = 1.137419834200614
$ws.Range("B25").Value = 0.920301176789792
$ws.Range("C25").Value = 0.2481865118657822
$ws.Range("E25").Value = 0.6220461076313839
$ws.Range("F25").Value = 1.765969994392776
$ws.Range("G25").Value = 0.2224555530345995
$ws.Range("H25").Value = 0.4032351895155983
$ws.Range("J25").Value = 0.02546209385334208
$ws.Range("M25").Value = 0.6467402572518921
$ws.Range("O25").Value = 1.164317125648708
